$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right column 4 -> 5, Wrong column -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right column 72 -> 90, Wrong column -1 -> -1.2, Max text updated
$ws.Range("B12").Value = 90
$ws.Range("C12").Value = -1.2
$ws.Range("E12").Value = "88.8/140"
